# Add a new "MetodoEntregaProducto" column (D) to the AgregarProducto
# parameters sheet, with a "Recoge_en_tienda" value row, matching the
# existing header/value formatting, then update the UI selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell D1 -------------------------------------------------
# Clone the formatting of C1 (bold header style) onto D1, then set text.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "MetodoEntregaProducto"

# --- Value cell D2 ----------------------------------------------------
# Clone the formatting of A2 (plain centered value style, not the
# hyperlink style used by C2) onto D2, then set text.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D2").Value = "Recoge_en_tienda"

# --- Column width -------------------------------------------------
# Widen column D to fit its new contents (matches the bestFit sizing
# used by the other parameter columns).
$ws.Range("D1:D2").ColumnWidth = 25

# --- Selection ------------------------------------------------------
$ws.Range("I6").Select()
